$wb = $excel.ActiveWorkbook

$wsUsers = $wb.Worksheets.Item("Users")
$wsCards = $wb.Worksheets.Item("Cards")

# --- Users sheet ("Sheet1") ---
# H8 changes from "admin" to "user"
$wsUsers.Range("H8").Value = "user"

# New row 10: user #9, only A/B/C/H populated
$wsUsers.Range("A10").Value = 9
$wsUsers.Range("B10").Value = "user"
$wsUsers.Range("C10").Value = "user"
$wsUsers.Range("H10").Value = "user"

# --- Cards sheet ("Sheet2") ---
# New column I: "user" enabled flag
$wsCards.Range("I2").Value = "enable"
$wsCards.Range("I3").Value = "yes"
$wsCards.Range("I4").Value = "yes"
$wsCards.Range("I5").Value = "yes"
$wsCards.Range("I6").Value = "yes"
$wsCards.Range("I7").Value = "yes"
$wsCards.Range("I8").Value = "yes"
$wsCards.Range("I9").Value = "no"
$wsCards.Range("I10").Value = "no"

# --- Selections / active sheet ---
# Final state: Users selection at B5 (not the active tab),
# Cards selection at K4 (Cards is the active/visible tab)
$wsUsers.Range("B5").Select() | Out-Null
$wsCards.Activate() | Out-Null
$wsCards.Range("K4").Select() | Out-Null
